$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row before row 13 (old rows 13-23 shift down to 14-24) ---
$ws.Rows.Item(13).EntireRow.Insert()

# The engine copies row-above formatting into the new row's A13 cell; reset it to
# a truly blank/default cell (matching the target, which has no A13 at all).
$ws.Range("D1").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").ClearContents()

# --- Row 10 ("Objetivos:"): replace the professor name with the real objectives text ---
$ws.Range("B10").Value2 = 'Conferir aos alunos uma visão geral da Indústria Siderúrgica, Metalúrgica, Mecânica e correlatas, bem como das principais características dos processos e arranjos produtivos destas indústrias.'
$ws.Range("C10").Value2 = 'Conferir aos alunos uma visão geral da Indústria Siderúrgica, Metalúrgica, Mecânica e correlatas, bem como das principais características dos processos e arranjos produtivos destas indústrias.'

# --- Row 13 (new, under "Docentes responsaveis:"): professor name moves here ---
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value2 = '5840560 - Marco Antonio Carvalho Pereira'
$ws.Range("C13").Value2 = '5840560 - Marco Antonio Carvalho Pereira'

# --- Row 14 (was row 13, "Programa resumido:"): new summary text + row height 60 ---
$ws.Range("B14").Value2 = '1) Processos da Indústria Siderúrgica. 2) Processos da Indústria Metalúrgica. 3) Processos da Indústria Mecânica.  4) Processos Industriais em Geral, exceto da Indústria Química'
$ws.Range("C14").Value2 = '1) Processos da Indústria Siderúrgica. 2) Processos da Indústria Metalúrgica. 3) Processos da Indústria Mecânica.  4) Processos Industriais em Geral, exceto da Indústria Química'
$ws.Rows.Item(14).RowHeight = 60

# --- Row 16 (was row 15, "Programa:"): new full program text ---
$ws.Range("B16").Value2 = '1) Processos da Indústria Siderúrgica: Obtenção de Aços. Demais Processos.2) Processos da Indústria Metalúrgica: Processos de Fundição. Demais Processos.3) Processos da Indústria Mecânica: Processos de Conformação Plástica. Conformação por Corte de Usinagem. Demais Processos4) Processos Industriais em geral: Indústria da Construção Civil Indústria farmacêutica, Indústria Automobilística, dentre outras.'
$ws.Range("C16").Value2 = '1) Processos da Indústria Siderúrgica: Obtenção de Aços. Demais Processos.2) Processos da Indústria Metalúrgica: Processos de Fundição. Demais Processos.3) Processos da Indústria Mecânica: Processos de Conformação Plástica. Conformação por Corte de Usinagem. Demais Processos4) Processos Industriais em geral: Indústria da Construção Civil Indústria farmacêutica, Indústria Automobilística, dentre outras.'

# --- Row 19 (was row 18, "Metodo:"): gains the teaching-method text ---
$ws.Range("B19").Value2 = 'Aulas expositivas. Seminários. Palestras feiras por profissionais de indústrias. Trabalhos em grupo. Debates e palestras.'
$ws.Range("C19").Value2 = 'Aulas expositivas. Seminários. Palestras feiras por profissionais de indústrias. Trabalhos em grupo. Debates e palestras.'

# --- Row 22 (was row 21, "Bibliografia:"): new bibliography text ---
$ws.Range("B22").Value2 = '1. Marcelo Breda Mourão et al. Introdução à Siderurgia, ABM, São Paulo, 20072. Fathi Habashi. Extractive Metallurgy, Gordon and Breach Science Publishers, 1986. 3. Luiz Antônio de Araújo. Manual de siderurgia - produção, Editora Arte & Ciência, São Paulo, 1997. 4. Alan H. Cottrell. Introdução à metalurgia, 2a edição, Fundação Calouste Gulbenkian, Lisboa, 1975.5. ASM Handbook Vol. 15 Casting - 1988 , Foundry Technology P.R. Beeley, 19726. John Campbell. Casting Butterworth-Heinemann, 19917. M. Siegel, Fundição. ABM, S.Paulo, 1979. 8. Amauri Garcia. Solidificação: Fundamentos e Aplicações, Editora da UNICAMP, Campinas, SP, 20089. Mauricio Prates de Campos Filho e Graeme John Davies Solidificação e Fundição de Metais e suas Ligas, Livros Técnicos e Científicos, Rio de Janeiro.10. AVITZUR, B. Metal Forming: processes and analysis – TATA Mc Graw-Hill Publishing Company Limited; New Delhi, 1977.11. RODRIGUES, J.; MARTINS, P. Tecnologia Mecânica: Tecnologia da deformação plástica. Aplicações Industriais. Escolar Editora, v.1 e v.2, 2010.12. CETLIN, P.R.; HELMAN, H. Fundamentos da conformação mecânicas dos metais. Ed. Artliber Ltda, 260p., 2005.13. BRESCIANI FILHO, E.; ZAVAGLIA, C.A.C.; NERY, F.A.C.; BOTTON, S.T. Conformação plástica dos metais. Ed. Unicamp, v.1 e v.2, 1986.14. DINIZ, A.E.; MARCONDES, F.C.; COPPINI, N.L. Tecnologia da usinagem dos materiais. Ed. Artlebet Ltda., 244p., 2000.'
$ws.Range("C22").Value2 = '1. Marcelo Breda Mourão et al. Introdução à Siderurgia, ABM, São Paulo, 20072. Fathi Habashi. Extractive Metallurgy, Gordon and Breach Science Publishers, 1986. 3. Luiz Antônio de Araújo. Manual de siderurgia - produção, Editora Arte & Ciência, São Paulo, 1997. 4. Alan H. Cottrell. Introdução à metalurgia, 2a edição, Fundação Calouste Gulbenkian, Lisboa, 1975.5. ASM Handbook Vol. 15 Casting - 1988 , Foundry Technology P.R. Beeley, 19726. John Campbell. Casting Butterworth-Heinemann, 19917. M. Siegel, Fundição. ABM, S.Paulo, 1979. 8. Amauri Garcia. Solidificação: Fundamentos e Aplicações, Editora da UNICAMP, Campinas, SP, 20089. Mauricio Prates de Campos Filho e Graeme John Davies Solidificação e Fundição de Metais e suas Ligas, Livros Técnicos e Científicos, Rio de Janeiro.10. AVITZUR, B. Metal Forming: processes and analysis – TATA Mc Graw-Hill Publishing Company Limited; New Delhi, 1977.11. RODRIGUES, J.; MARTINS, P. Tecnologia Mecânica: Tecnologia da deformação plástica. Aplicações Industriais. Escolar Editora, v.1 e v.2, 2010.12. CETLIN, P.R.; HELMAN, H. Fundamentos da conformação mecânicas dos metais. Ed. Artliber Ltda, 260p., 2005.13. BRESCIANI FILHO, E.; ZAVAGLIA, C.A.C.; NERY, F.A.C.; BOTTON, S.T. Conformação plástica dos metais. Ed. Unicamp, v.1 e v.2, 1986.14. DINIZ, A.E.; MARCONDES, F.C.; COPPINI, N.L. Tecnologia da usinagem dos materiais. Ed. Artlebet Ltda., 244p., 2000.'
